$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5; this shifts rows 5-30 down to 6-31,
# preserving all of their existing values (matching the diff's row shift).
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new data record.
$ws.Cells.Item(5, 1).Value = 10
$ws.Cells.Item(5, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(5, 3).Value = "La Araucanía"
$ws.Cells.Item(5, 4).Value = 45163
$ws.Cells.Item(5, 5).Value = 9
$ws.Cells.Item(5, 6).Value = 100112017
$ws.Cells.Item(5, 7).Value = "Ramas de apio"
$ws.Cells.Item(5, 8).Value = "Sin especificar"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 30
$ws.Cells.Item(5, 11).Value = 5000
$ws.Cells.Item(5, 12).Value = 5000
$ws.Cells.Item(5, 13).Value = 5000
$ws.Cells.Item(5, 14).Value = "$/paquete"
$ws.Cells.Item(5, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(5, 16).Value = 5000
$ws.Cells.Item(5, 17).Value = 1
$ws.Cells.Item(5, 18).Value = "Hortaliza"
